$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Add Devices On Ethernet")
$ws2 = $wb.Worksheets.Item("Add Devices on RBus")
$ws3 = $wb.Worksheets.Item("Sheet1")

# Relocate the last data row (row 11, the "Pro32xD" entry) of "Add Devices On
# Ethernet" onto "Sheet1" at row 5 (values + formatting), then tag it with a
# new note cell pointing back at the source sheet.
$ws1.Range("A11:I11").Copy($ws3.Range("A5"))
$ws3.Cells.Item(5, 11).Value = "Add Devices On Ethernet"

# "Add Devices On Ethernet" row 10 now documents the Pro32xD panel instead of
# the Pro16xD one.
$ws1.Cells.Item(10, 1).Value = "Pro32xD"

# That row has now been relocated, so drop it from "Add Devices On Ethernet".
$ws1.Rows.Item(11).Delete()

# The extra duplicate data row on "Add Devices on RBus" is no longer needed.
$ws2.Rows.Item(11).Delete()

# Restore each sheet's cursor/selection and make "Add Devices On Ethernet" the
# active tab again.
$ws2.Range("A10").Select()
$ws3.Range("K5").Select()

$ws1.Activate()
$ws1.Range("A8").Select()
